# Update the daily COVID-19 case-count rows (542-552) on the single sheet.
# Columns B, H, J, K are volatile TODAY()-driven shared formulas that
# recompute automatically once the inputs below are written, so only the
# literal input cells (C, E, F, G and, for the newly-populated rows, L/M)
# need to be touched here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L (12) and M (13) are formatted as Text (@) in this sheet, so a
# plain numeric assignment through COM would land as a text "0" (matching
# real Excel's "typing a number into a text cell" behaviour). Flip the
# number format to General for the instant of the write, then restore the
# original Text format, so the stored value stays a genuine number (as the
# target file expects).
function Set-NumericValue($cell, $value) {
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "General"
    $cell.Value2 = $value
    $cell.NumberFormat = $fmt
}

# Row 542
$ws.Cells.Item(542, 3).Value2 = 91

# Row 544
$ws.Cells.Item(544, 5).Value2 = 5
$ws.Cells.Item(544, 7).Value2 = 15

# Row 545
$ws.Cells.Item(545, 3).Value2 = 117
$ws.Cells.Item(545, 5).Value2 = 5
$ws.Cells.Item(545, 7).Value2 = 16

# Row 546
$ws.Cells.Item(546, 3).Value2 = 92
$ws.Cells.Item(546, 5).Value2 = 7
$ws.Cells.Item(546, 7).Value2 = 19

# Row 547
$ws.Cells.Item(547, 3).Value2 = 131
$ws.Cells.Item(547, 5).Value2 = 7
$ws.Cells.Item(547, 7).Value2 = 23

# Row 548
$ws.Cells.Item(548, 3).Value2 = 85
$ws.Cells.Item(548, 5).Value2 = 6
$ws.Cells.Item(548, 7).Value2 = 20

# Row 549 (previously-blank trailing row, now gets real daily figures)
$ws.Cells.Item(549, 3).Value2 = 81
$ws.Cells.Item(549, 5).Value2 = 5
$ws.Cells.Item(549, 6).Value2 = 4
$ws.Cells.Item(549, 7).Value2 = 16
Set-NumericValue $ws.Cells.Item(549, 12) 0
Set-NumericValue $ws.Cells.Item(549, 13) 0

# Row 550
$ws.Cells.Item(550, 3).Value2 = 41
$ws.Cells.Item(550, 5).Value2 = 8
$ws.Cells.Item(550, 6).Value2 = 5
$ws.Cells.Item(550, 7).Value2 = 16
Set-NumericValue $ws.Cells.Item(550, 12) 0
Set-NumericValue $ws.Cells.Item(550, 13) 0

# Row 551
$ws.Cells.Item(551, 3).Value2 = 31
$ws.Cells.Item(551, 5).Value2 = 8
$ws.Cells.Item(551, 6).Value2 = 5
$ws.Cells.Item(551, 7).Value2 = 17
Set-NumericValue $ws.Cells.Item(551, 12) 0
Set-NumericValue $ws.Cells.Item(551, 13) 0

# Row 552
$ws.Cells.Item(552, 3).Value2 = 9
$ws.Cells.Item(552, 5).Value2 = 8
$ws.Cells.Item(552, 6).Value2 = 5
$ws.Cells.Item(552, 7).Value2 = 20
Set-NumericValue $ws.Cells.Item(552, 12) 0
Set-NumericValue $ws.Cells.Item(552, 13) 0

# Window state: the sheet is frozen at B3 (1 column / 2 rows). The author
# scrolled the bottom-right pane down and left the selection on F546;
# reproduce the selection (the one piece of view state this COM surface
# can actually move) without disturbing the freeze itself.
$ws.Activate()
$target = $ws.Range("F546")
$excel.Goto($target, $true)
$target.Select()
